$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the second data row (MkPV / Mouse kidney parvoviruses entry) first,
# so the remaining edits below intern new shared strings in the same order
# Excel used (Ichthama virus details before the genus name).
$ws.Rows(3).Delete()

# Row 2 now becomes the single remaining data row. Refactor it from the
# Chaphamaparvovirus (Chicken chapparvovirus-2 / CChPV) entry into the new
# Icthamaparvovirus (Ichthyic parvovirus isolate HMU-HKU / IcthPV) entry.
$ws.Range("A2").Value = "MN162688"
$ws.Range("B2").Value = "Hamaparvovirinae"
$ws.Range("D2").Value = "IcthPV "
$ws.Range("E2").Value = "Ichthyic parvovirus isolate HMU-HKU"
$ws.Range("C2").Value = "Icthamaparvovirus"
$ws.Range("F2").Value = "n/a"
$ws.Range("G2").Value = "Tilapia"

# Move the active selection to match the author's final cursor position.
$ws.Range("G20").Select()
